# Sheet update for Ist year 2k22
# Adds a new worksheet "sheet5" listing the First Year 2k22 volunteer members.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end of the workbook ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "sheet5"

$names = @(
    "Abhinamyu Singh",
    "Aditya Sharma",
    "Ankur Yadav",
    "Anshuman Payasi",
    "Archan Banerjee",
    "Ayan Chadoria",
    "Ayushi Shukla",
    "Gugli Thakur",
    "Kunal Dhiman",
    "Laksh Bhandari",
    "Mehul Ambastha",
    "Naman Sharma",
    "Oshin Sharma",
    "Prikshit Saini",
    "Rishika Sharma",
    "Sakshi Gothwal",
    "Sana Sheikh",
    "Tanashvi Joshi",
    "Urvashi Pandey",
    "Vishesh Garg"
)

# --- Header row (row 1) -----------------------------------------------------
# Same headers as the other sheets except column F, which is named
# "img_link" here instead of "image_link".
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "year"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "link"
$ws.Range("E1").Value = "img"
$ws.Range("F1").Value = "img_link"
$ws.Range("G1").Value = "branch"
$ws.Range("H1").Value = "city"
$ws.Range("I1").Value = "state"
$ws.Range("J1").Value = "about"
$ws.Rows.Item(1).RowHeight = 15

# --- Member rows (rows 2-21) -------------------------------------------------
# Row 2 first (name + year), then the remaining names down column A, then the
# remaining years down column B, then the description down column C - this
# mirrors the member sheet's original authoring order.
$ws.Range("A2").Value = $names[0]
$ws.Range("B2").Value = "First Year"

for ($i = 1; $i -lt $names.Count; $i++) {
    $ws.Range("A" + ($i + 2)).Value = $names[$i]
}

for ($i = 1; $i -lt $names.Count; $i++) {
    $ws.Range("B" + ($i + 2)).Value = "First Year"
}

for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Range("C" + ($i + 2)).Value = "Volunteer Member"
}

# Style rows 3-21 (Arial 10, left aligned) and give them a taller row height,
# matching the formatting used for the rest of the member table.
for ($r = 3; $r -le 21; $r++) {
    $rng = $ws.Range("A" + $r + ":C" + $r)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.HorizontalAlignment = -4131
    $ws.Rows.Item($r).RowHeight = 18.75
}

# --- Trailing styled-but-empty rows (22-28), column C only -----------------
for ($r = 22; $r -le 28; $r++) {
    $c = $ws.Range("C" + $r)
    $c.Font.Name = "Arial"
    $c.Font.Size = 10
    $c.HorizontalAlignment = -4131
    $ws.Rows.Item($r).RowHeight = 18.75
}

# --- View / selection state ---------------------------------------------
$ws.Activate()
$ws.Range("C21").Select()
